$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33
$ws.Range("B3").Value = 33
$ws.Range("B4").Value = 33

$ws.Rows(5).Delete()
